$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1, copying the formatting of C1 (header style) and
# then setting its text to the new "Category" prompt.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Category (names OR ids)"

# Add new data cell D2, copying the formatting of C2 (blank input cell style).
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Widen the new column D to fit the longer header text.
$ws.Columns.Item(4).ColumnWidth = 18.52
